$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Commit the previous generator row (Legal_UPD) output into history rows 9:11
# with the new "Hotline" NYS Smokers Quit Line entry.
$ws.Range("A9").Value = "mod_Accordion_ui('Hotline_NYSQuitline')"
$ws.Range("A10").Value = "mod_Accordion_server('Hotline_NYSQuitline', selector=selection, data=HLdata(), title = c('NYS Smokers Quit Line'), Visible = T)"
$ws.Range("A11").Value = "mod_info_server('Hotline_NYSQuitline', selector = selection, data = HLdata(), rownametitle = c('NY Quits - Smokers Quit Line'), phone = T, website = T)"

# Update the generator input row (row 2) with the new AOD entry.
$ws.Range("A2").Value = "AOD_NYSQuitline"
$ws.Range("B2").Value = "NY Quits - Smokers Quit Line"
$ws.Range("C2").Value = "AODdata()"

# D2 switches from the hyperlink-like style it had to the plain Times New
# Roman style already used by A2/B2/E2/F2, so copy that formatting over.
$ws.Range("A2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D2").Value = "NYS Smokers Quit Line"

# Commit the newly generated AOD output into the next free history rows 14:16.
$ws.Range("A14").Value = "mod_Accordion_ui('AOD_NYSQuitline')"
$ws.Range("A15").Value = "mod_Accordion_server('AOD_NYSQuitline', selector=selection, data=AODdata(), title = c('NYS Smokers Quit Line'), Visible = T)"
$ws.Range("A16").Value = "mod_info_server('AOD_NYSQuitline', selector = selection, data = AODdata(), rownametitle = c('NY Quits - Smokers Quit Line'), phone = T, website = T)"

# Match the new cell font (Times New Roman 12pt, matching style index 3) for A14:A16
$ws.Range("A14:A16").Font.Name = $ws.Range("A9").Font.Name
$ws.Range("A14:A16").Font.Size = $ws.Range("A9").Font.Size

# Update selection to reflect the newly active history rows
$excel.Goto($ws.Range("A14:A16"))
